# Apply updated crypto market data (prices & 1h volume change) per Wed Jun 26 2024 refresh
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '61.085.26'
$ws.Range('E2').Value = '  -1.32%  '
$ws.Range('D3').Value = '3.366.64'
$ws.Range('E3').Value = '  -1.32%  '
$ws.Range('E4').Value = '  +0.14%  '
$ws.Range('D5').Value = '''570.73'
$ws.Range('E5').Value = '  -1.14%  '
$ws.Range('D6').Value = '''135.74'
$ws.Range('E6').Value = '  -1.92%  '
$ws.Range('E7').Value = '  +0.13%  '
$ws.Range('D8').Value = '3.364.68'
$ws.Range('E8').Value = '  -1.42%  '
$ws.Range('D9').Value = '''0.469'
$ws.Range('E9').Value = '  -1.91%  '
$ws.Range('D10').Value = '''7.50'
$ws.Range('E10').Value = '  -0.24%  '
$ws.Range('E11').Value = '  -3.77%  '
$ws.Range('E12').Value = '  -3.00%  '
$ws.Range('D13').Value = '3.936.67'
$ws.Range('E13').Value = '  -1.35%  '
$ws.Range('E14').Value = '  +1.24%  '
$ws.Range('D15').Value = '''25.97'
$ws.Range('E15').Value = '  +1.97%  '
$ws.Range('E16').Value = '  -4.61%  '
$ws.Range('D17').Value = '3.359.02'
$ws.Range('E17').Value = '  -1.46%  '
$ws.Range('D18').Value = '61.203.91'
$ws.Range('E18').Value = '  -0.97%  '
$ws.Range('D19').Value = '''14.00'
$ws.Range('E19').Value = '  -1.01%  '
$ws.Range('E20').Value = '  -1.91%  '
$ws.Range('D21').Value = '''9.24'
$ws.Range('E21').Value = '  -2.28%  '
$ws.Range('D22').Value = '''376.71'
$ws.Range('E22').Value = '  -3.90%  '
$ws.Range('D23').Value = '''0.553'
$ws.Range('E23').Value = '  -3.57%  '
$ws.Range('D24').Value = '3.493.42'
$ws.Range('E24').Value = '  -1.47%  '
$ws.Range('D25').Value = '''1.00'
$ws.Range('E25').Value = '  -0.18%  '
$ws.Range('B26').Value = 'PEPE'
$ws.Range('C26').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D26').Value = '''0.0000125'
$ws.Range('E26').Value = '  -3.27%  '
$ws.Range('B27').Value = 'Litecoin'
$ws.Range('C27').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('D27').Value = '''71.04'
$ws.Range('E27').Value = '  -0.59%  '
$ws.Range('D28').Value = '''1.76'
$ws.Range('E28').Value = '  +9.79%  '
$ws.Range('E29').Value = '  +0.75%  '
$ws.Range('D30').Value = '''7.46'
$ws.Range('E30').Value = '  -3.27%  '
$ws.Range('D31').Value = '''0.166'
$ws.Range('E31').Value = '  +4.16%  '
$ws.Range('D32').Value = '''8.13'
$ws.Range('E32').Value = '  -2.04%  '
$ws.Range('E33').Value = '  -1.56%  '
$ws.Range('D35').Value = '''23.53'
$ws.Range('E35').Value = '  -0.27%  '
$ws.Range('D36').Value = '''5.19'
$ws.Range('E36').Value = '  -6.09%  '
$ws.Range('D37').Value = '''6.76'
$ws.Range('E37').Value = '  -3.14%  '
$ws.Range('D38').Value = '''164.90'
$ws.Range('E38').Value = '  +2.13%  '
$ws.Range('E39').Value = '  -2.19%  '
$ws.Range('E40').Value = '  -5.34%  '
$ws.Range('E41').Value = '  +0.18%  '
$ws.Range('D42').Value = '''0.769'
$ws.Range('E42').Value = '  -0.80%  '
$ws.Range('E43').Value = '  -2.34%  '
$ws.Range('D44').Value = '''41.47'
$ws.Range('E44').Value = '  -0.17%  '
$ws.Range('B45').Value = 'Filecoin'
$ws.Range('C45').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D45').Value = '''4.39'
$ws.Range('E45').Value = '  -2.01%  '
$ws.Range('B46').Value = 'ONDO'
$ws.Range('C46').Value = 'https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo'
$ws.Range('D46').Value = '''1.20'
$ws.Range('E46').Value = '  -2.23%  '
$ws.Range('D47').Value = '''23.86'
$ws.Range('E47').Value = '  -5.30%  '
$ws.Range('D48').Value = '''23.23'
$ws.Range('E48').Value = '  +1.35%  '
$ws.Range('D49').Value = '''6.80'
$ws.Range('E49').Value = '  -2.67%  '
$ws.Range('D50').Value = '2.343.35'
$ws.Range('E50').Value = '  -1.36%  '
$ws.Range('B51').Value = 'LidoDAOToken'
$ws.Range('C51').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D51').Value = '''2.35'
$ws.Range('E51').Value = '  +1.59%  '
